# docgen.hangul.auto_postfix -> docgen.hangul.auto_suffix 로 수정
$wb = $excel.ActiveWorkbook

$wsEn = $wb.Worksheets.Item("en")
$wsKo = $wb.Worksheets.Item("ko")

# Rename the variable name in column A, row 13, on both sheets.
$wsEn.Range("A13").Value = "docgen.hangul.auto_suffix"
$wsKo.Range("A13").Value = "docgen.hangul.auto_suffix"

# Update the selection on each sheet to A14, and switch the active sheet
# from "ko" to "en".
$wsKo.Range("A14").Select()
$wsEn.Activate()
$wsEn.Range("A14").Select()
